$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E (shifts old batsman..sr columns to F..K)
$ws.Range("D1:E1").EntireColumn.Insert()

# New header cells for the inserted columns
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Fill in the new ownTeam/oppTeam values for the existing data row (row 2)
$ws.Range("D2").Value = "Rajasthan Royals"
$ws.Range("E2").Value = "Royal Challengers Bangalore"

# Ensure numeric-looking text values stay as text (matches t="str" source data)
$ws.Range("G2:K2").NumberFormat = "@"
$ws.Range("G3:K4").NumberFormat = "@"

# New row 3 data
$ws.Range("A3").Value = " Sharjah"
$ws.Range("B3").Value = " October 09 2020"
$ws.Range("C3").Value = "Capitals won by 46 runs"
$ws.Range("D3").Value = "Rajasthan Royals"
$ws.Range("E3").Value = "Delhi Capitals"
$ws.Range("F3").Value = "Mahipal Lomror "
$ws.Range("G3").Value = "1"
$ws.Range("H3").Value = "2"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "50.00"

# New row 4 data
$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " October 06 2020"
$ws.Range("C4").Value = "Mumbai won by 57 runs"
$ws.Range("D4").Value = "Rajasthan Royals"
$ws.Range("E4").Value = "Mumbai Indians"
$ws.Range("F4").Value = "Mahipal Lomror "
$ws.Range("G4").Value = "11"
$ws.Range("H4").Value = "13"
$ws.Range("I4").Value = "1"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "84.61"
